$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated metadata obtained 2016-04-06:
#  - B3: measure -> dimension for the "sector-vab-descripcion" concept
#  - B4: type column changed from "medida" (measure) to "dim" (dimension)
$ws.Range("B3").Value = "iaest-dimension:sector-vab-descripcion"
$ws.Range("B4").Value = "dim"

# New row describing the mapping file for this dimension's codes.
# Copy B5's formatting (style index 1, Arial 10) onto the new B6 cell
# before writing its value, so the new cell matches the sheet's styling.
$ws.Range("B5").Copy($ws.Range("B6"))
$ws.Range("B6").Value = "mapping-sector-vab-descripcion.xlsx"
